$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($sheet, $addr, $val) {
    $sheet.Range($addr).NumberFormat = "@"
    $sheet.Range($addr).Value = $val
    $sheet.Range($addr).ClearFormats()
}

$ws.Range("D2").Value = "26.804.98"
$ws.Range("E2").Value = "  -0.07%  "
$ws.Range("D3").Value = "1.638.89"
$ws.Range("E3").Value = "  -0.49%  "
$ws.Range("E4").Value = "  -0.66%  "
Set-TextValue $ws "D5" "218.96"
$ws.Range("E5").Value = "  +0.72%  "
Set-TextValue $ws "D6" "0.500"
$ws.Range("E6").Value = "  -0.69%  "
$ws.Range("E7").Value = "  -0.61%  "
$ws.Range("E8").Value = "  -0.53%  "
Set-TextValue $ws "D9" "0.0622"
$ws.Range("E9").Value = "  -0.87%  "
Set-TextValue $ws "D10" "19.26"
$ws.Range("E10").Value = "  +0.14%  "
$ws.Range("E11").Value = "  +0.22%  "
$ws.Range("D12").Value = "1.865.89"
$ws.Range("E12").Value = "  -0.51%  "
$ws.Range("D13").Value = "1.635.13"
$ws.Range("E13").Value = "  -0.53%  "
Set-TextValue $ws "D14" "4.14"
$ws.Range("E14").Value = "  -1.19%  "
Set-TextValue $ws "D15" "0.525"
$ws.Range("E15").Value = "  -0.65%  "
Set-TextValue $ws "D16" "64.72"
$ws.Range("E16").Value = "  -0.06%  "
$ws.Range("D17").Value = "26.789.67"
$ws.Range("E17").Value = "  -0.10%  "
$ws.Range("D18").Value = "0.0₃0734"
$ws.Range("E18").Value = "  -0.57%  "
Set-TextValue $ws "D19" "214.82"
$ws.Range("E19").Value = "  +0.13%  "
$ws.Range("E20").Value = "  -0.54%  "
$ws.Range("E21").Value = "  -0.02%  "
Set-TextValue $ws "D22" "6.30"
$ws.Range("E22").Value = "  +0.03%  "
Set-TextValue $ws "D23" "2.36"
$ws.Range("E23").Value = "  -2.97%  "
Set-TextValue $ws "D24" "9.11"
$ws.Range("E24").Value = "  -2.86%  "
Set-TextValue $ws "D25" "147.46"
$ws.Range("E25").Value = "  +1.63%  "
$ws.Range("E26").Value = "  -0.76%  "
Set-TextValue $ws "D27" "0.119"
$ws.Range("E27").Value = "  +0.07%  "
$ws.Range("E28").Value = "  -0.84%  "
$ws.Range("E29").Value = "  -0.22%  "
$ws.Range("E30").Value = "  -1.48%  "
$ws.Range("E31").Value = "  +1.22%  "
Set-TextValue $ws "D32" "3.39"
$ws.Range("E32").Value = "  +2.00%  "
$ws.Range("E33").Value = "  -0.01%  "
$ws.Range("E34").Value = "  +0.11%  "
$ws.Range("D35").Value = "1.262.58"
$ws.Range("E35").Value = "  -1.78%  "
$ws.Range("E37").Value = "  -0.21%  "
$ws.Range("E38").Value = "  -1.91%  "
Set-TextValue $ws "D39" "0.813"
$ws.Range("E39").Value = "  -1.85%  "
$ws.Range("E40").Value = "  -0.50%  "
Set-TextValue $ws "D41" "0.804"
$ws.Range("E41").Value = "  -1.18%  "
$ws.Range("E42").Value = "  -0.61%  "
$ws.Range("D43").Value = "1.777.39"
$ws.Range("E43").Value = "  -1.06%  "
$ws.Range("E44").Value = "  -4.58%  "
Set-TextValue $ws "D45" "92.13"
$ws.Range("E45").Value = "  +0.69%  "
Set-TextValue $ws "D46" "60.03"
$ws.Range("E46").Value = "  +0.52%  "
Set-TextValue $ws "D47" "1.59"
$ws.Range("E47").Value = "  -1.21%  "
$ws.Range("E48").Value = "  -0.56%  "
$ws.Range("B49").Value = "EnergySwap"
$ws.Range("C49").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
Set-TextValue $ws "D49" "7.53"
$ws.Range("E49").Value = "  -1.70%  "
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
Set-TextValue $ws "D50" "0.0963"
$ws.Range("E50").Value = "  -1.21%  "

Write-Output "edit complete"
